$d = $word.ActiveDocument

# 1) Expand/replace the sentence about antibiotic resistance in the
#    southeastern United States with the new speculative explanation.
$old = "Continued through our analysis, for the United States, antibiotic resistance appears to be greater among the southeastern portion of the country. On a global level, points of interest include Russia, India, the Middle East, and parts of Africa."
$new = "Continued through our analysis, for the United States, antibiotic resistance appears to be lower among the southeastern portion of the country. We speculate that the reason for this is because the number of prescribed antibiotics in the southeastern states is higher than any other part of the country which allows a variety of antibiotics to be distributed and counteract different bacteria. If we look at the world map, we can confirm this speculation since the percent of resistance is extremely low in the United States. Now on a global level, points of interest include Russia, India, the Middle East, and parts of Africa."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2) Remove 5 of the trailing empty paragraphs (leftover blank lines with
#    underline formatting) that used to pad out the space after the
#    "future research" paragraph near the end of the document.
$target = "There are several rich domains for future research"
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains($target)) {
        $anchorIndex = $i
        break
    }
}

$firstBlank = $anchorIndex + 1
for ($i = 1; $i -le 5; $i++) {
    $d.Paragraphs($firstBlank).Range.Delete()
}
